$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date Submitted value moves from Dec 4th to Dec 5th (serial 43073 -> 43074)
$ws.Range("B4").Value = 43074

# Fill in the outcome vector + actual estimates for File 5 through File 8
# (rows 42-45): D = 0.88 (estimate for the true outcome), E:P = 0.01 each
# (the other 12 possibilities) so each row sums to 1.
$rows = 42, 43, 44, 45
foreach ($r in $rows) {
    $ws.Range("D$r").Value = 0.88
    $ws.Range("E$r`:P$r").Value = 0.01
}

# Move the view back to the top and select B5 instead of B49
$ws.Range("B5").Select()
